$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Title: "Recommendation Model: FIFA Players" (merge runs -> identical text,
#    nothing to actually change content-wise; skip).
# ---------------------------------------------------------------------------

# ---------------------------------------------------------------------------
# 2) Overview / Introduction paragraph: extend with KNN sentence.
# ---------------------------------------------------------------------------
$old1 = "First of all, data science nowadays is becoming more important than ever. It the way to an artificial intelligence world while using machine learning techniques to automate, predict and solve problems. In this project, I will be building a recommendation system. There are two most important models of recommendation, which are Content-based recommendation and collaborative-based recommendation."
$new1 = "First of all, data science nowadays is becoming more important than ever. It the way to an artificial intelligence world while using machine learning techniques to automate, predict and solve problems. In this project, I will be building a recommendation system by using the K Nearest Neighbor (KNN) algorithm."
$d.Content.Find.Execute($old1, $false, $false, $false, $false, $false, $true, 1, $false, $new1, 2) | Out-Null

# ---------------------------------------------------------------------------
# 3) Question / Need paragraph: rewrite ending about KNN usage.
# ---------------------------------------------------------------------------
$old2 = "My project is building a recommendation system that recommend multiple players, based on the player being search by. In this project, I will be using their similarities to recommend. "
$new2 = "My project is building a recommendation system that recommend multiple players, based on the player being search by. In this project, I will be using knn to their nearest players based on their similarities. "
$d.Content.Find.Execute($old2, $false, $false, $false, $false, $false, $true, 1, $false, $new2, 2) | Out-Null

# ---------------------------------------------------------------------------
# 4) Tools section: insert a new "Seaborn: Visualization" paragraph right
#    before the "Scikit Learn: Modeling" paragraph.
# ---------------------------------------------------------------------------
$found = $false
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -eq "Scikit Learn: Modeling`r") {
        $p.Range.InsertParagraphBefore()
        $newPara = $d.Paragraphs.Item($i)
        $newPara.Range.Text = "Seaborn: Visualization"
        $found = $true
        break
    }
}
if (-not $found) {
    Write-Host "WARNING: Scikit Learn paragraph not found"
}
